$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, alignment) from H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    5,6
    11,11
    7,8
    2,2
    7,7
    5,5
    7,7
    9,9
    10,10
    7,7
    8,8
    7,8
    8,8
    6,6
    8,8
    8,8
    7,7
    8,8
    8,8
    8,8
    8,8
    7,7
    3,3
    4,4
)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i * 2]
    $ws.Cells.Item($row, 10).Value = $data[$i * 2 + 1]
}
